# Inserts a new "id" column at the front of the inventory table (producto /
# cantidad / ultima_actualizacion shift one column to the right) and updates
# the "cantidad" / "ultima_actualizacion" values to the refreshed dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shift producto/cantidad/ultima_actualizacion one column to the right by
#    inserting a fresh column A.
$ws.Columns.Item(1).Insert()

# 2. New header + id values (0-based row index) for column A. Match the
#    bold/bordered/centered look of the other header cells (B1/C1/D1) by
#    copying their format onto the new header cell.
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 1).Value = "id"

$ids = 0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17
for ($i = 0; $i -lt $ids.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $ids[$i]
}

# 3. Updated cantidad (col C) / ultima_actualizacion (col D) values. The date
#    strings must stay literal text (not get auto-converted to date serials),
#    so format the cell as text before writing, then restore a plain/default
#    style so no stray number format sticks around on the cell.
$data = @(
    @{ Row = 2;  Cantidad = 187; Fecha = "2025-06-18" },
    @{ Row = 3;  Cantidad = 113; Fecha = "2025-06-25" },
    @{ Row = 4;  Cantidad = 114; Fecha = "2025-06-26" },
    @{ Row = 5;  Cantidad = 244; Fecha = "2025-06-25" },
    @{ Row = 6;  Cantidad = 184; Fecha = "2025-06-20" },
    @{ Row = 7;  Cantidad = 98;  Fecha = "2025-06-25" },
    @{ Row = 8;  Cantidad = 154; Fecha = "2025-06-22" },
    @{ Row = 9;  Cantidad = 96;  Fecha = "2025-06-15" },
    @{ Row = 10; Cantidad = 140; Fecha = "2025-06-14" },
    @{ Row = 11; Cantidad = 140; Fecha = "2025-06-26" },
    @{ Row = 12; Cantidad = 159; Fecha = "2025-06-18" },
    @{ Row = 13; Cantidad = 109; Fecha = "2025-06-26" },
    @{ Row = 14; Cantidad = 128; Fecha = "2025-06-26" },
    @{ Row = 15; Cantidad = 85;  Fecha = "2025-06-25" },
    @{ Row = 16; Cantidad = 101; Fecha = "2025-06-26" },
    @{ Row = 17; Cantidad = 84;  Fecha = "2025-06-21" },
    @{ Row = 18; Cantidad = 92;  Fecha = "2025-06-22" },
    @{ Row = 19; Cantidad = 80;  Fecha = "2025-06-26" }
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 3).Value = $entry.Cantidad

    $dateCell = $ws.Cells.Item($entry.Row, 4)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $entry.Fecha
    $dateCell.Style = "Normal"
}
